$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

for ($row = 3; $row -le 23; $row++) {
    $hCell = $ws.Cells.Item($row, 8)   # Column H - PERIOD TO EXPIRE
    $hVal = $hCell.Value()
    $hCell.Value = $hVal - 1

    $iCell = $ws.Cells.Item($row, 9)   # Column I - LAST UPDATE
    $iCell.NumberFormat = "@"
    $iCell.Value = "04-Nov-2025"
}
